$wb = $excel.ActiveWorkbook

# The source sheet used as a structural template for the two new sheets -
# it already carries the right column widths, header styling, page setup,
# etc. that the new "XXL" sheets were cloned from.
$template = $wb.Worksheets.Item($wb.Worksheets.Count)

# --- Create the two new sheets (in creation order: "(Round 2)" first, then
# "(2)") by copying the template sheet, then rename them. ---
$template.Copy($null, $template)
$wb.Worksheets.Item($wb.Worksheets.Count).Name = "(XXL) 460x256 (Round 2)"

$template.Copy($null, $wb.Worksheets.Item("(XXL) 460x256 (Round 2)"))
$wb.Worksheets.Item($wb.Worksheets.Count).Name = "(XXL) 460x256 (2)"

# Final tab order has "(2)" before "(Round 2)" -- move it into place.
# (NOTE: re-fetch sheet handles by name afterwards -- Move() invalidates
# previously-held references for the sheets whose positions shifted.)
$wb.Worksheets.Item("(XXL) 460x256 (2)").Move($wb.Worksheets.Item("(XXL) 460x256 (Round 2)"))

$xxl2 = $wb.Worksheets.Item("(XXL) 460x256 (2)")
$xxlRound2 = $wb.Worksheets.Item("(XXL) 460x256 (Round 2)")

# --- Fill in "(XXL) 460x256 (2)" ---
$xxl2.Range("A2:K20").ClearContents()

$xxl2.Range("J5").Value = "Keygen Offset (seconds):"
$xxl2.Range("K5").Value = 249

$xxl2.Range("A6").Value = 128
$xxl2.Range("B6").Value = 2
$xxl2.Range("C6").Value = 11.13
$xxl2.Range("D6").Formula = "=C6-((B6*K5)/60)"

$xxl2.Range("A7").Value = 128
$xxl2.Range("B7").Value = 0
$xxl2.Range("C7").Value = 1
$xxl2.Range("D7").Formula = "=C7-((B7*K5)/60)"

$xxl2.Range("A8").Value = 128
$xxl2.Range("B8").Value = 0
$xxl2.Range("C8").Value = 0.13
$xxl2.Range("D8").Formula = "=C8-((B8*K5)/60)"

$xxl2.Range("A9").Value = 128
$xxl2.Range("B9").Value = 0
$xxl2.Range("C9").Value = 0.36
$xxl2.Range("D9").Formula = "=C9-((B9*K5)/60)"
$xxl2.Range("F9").Formula = "=AVERAGE(C6:C9)"

$xxl2.Range("B12").Value = "Avgs:"
$xxl2.Range("C12").Formula = "=AVERAGE(C6:C9)"
$xxl2.Range("D12").Formula = "=AVERAGE(D6:D9)"
$xxl2.Range("F12").Formula = "=AVERAGE(F2:F9)"

$xxl2.Range("C3").Select()

# --- Fill in "(XXL) 460x256 (Round 2)" ---
$xxlRound2.Range("A2:K20").ClearContents()

$xxlRound2.Range("J5").Value = "Keygen Offset (seconds):"
$xxlRound2.Range("K5").Value = 249

$xxlRound2.Range("A6").Value = 128
$xxlRound2.Range("B6").Value = 0
$xxlRound2.Range("C6").Value = 0.33
$xxlRound2.Range("D6").Formula = "=C6-((B6*K5)/60)"

$xxlRound2.Range("A7").Value = 128
$xxlRound2.Range("B7").Value = 0
$xxlRound2.Range("C7").Value = 0.11
$xxlRound2.Range("D7").Formula = "=C7-((B7*K5)/60)"

$xxlRound2.Range("A8").Value = 128
$xxlRound2.Range("B8").Value = 0
$xxlRound2.Range("C8").Value = 0.45
$xxlRound2.Range("D8").Formula = "=C8-((B8*K5)/60)"

$xxlRound2.Range("A9").Value = 128
$xxlRound2.Range("B9").Value = 0
$xxlRound2.Range("C9").Value = 0.47
$xxlRound2.Range("D9").Formula = "=C9-((B9*K5)/60)"
$xxlRound2.Range("F9").Formula = "=AVERAGE(C6:C9)"

$xxlRound2.Range("B12").Value = "Avgs:"
$xxlRound2.Range("C12").Formula = "=AVERAGE(C6:C9)"
$xxlRound2.Range("D12").Formula = "=AVERAGE(D6:D9)"
$xxlRound2.Range("F12").Formula = "=AVERAGE(F2:F9)"

# Select last so this sheet ends up the active / tab-selected one.
$xxlRound2.Range("B7").Select()
